$d = $word.ActiveDocument

$old2 = "To offer a user-friendly and engaging digital learning platform that helps students grasp academic concepts with ease, prepare effectively for exams, and study at their own pace. Through high-quality video content, interactive quizzes, personalized learning paths, and instant academic assistance, we aim to support learners in building confidence and reaching their full academic potential" + [char]0x2014 + "anytime, from anywhere."
$new2 = "To offer a user-friendly and engaging digital learning platform's fully functional mock version created by sophomore students of Web programming course that helps students in their academic life."

$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null
